$wb = $excel.ActiveWorkbook

# --- 1. reference_controls sheet: column A (ref_id), rows 2..145
#     M<n> -> 5G-M<n>
$refSheet = $wb.Worksheets.Item("reference_controls")
$refUsed = $refSheet.UsedRange
$refRowCount = $refUsed.Rows.Count

for ($r = 2; $r -le $refRowCount; $r++) {
    $cell = $refSheet.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = $val -replace '^M(\d+)$', '5G-M$1'
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# --- 2. requirements sheet: column F (reference_controls), every row
#     Each cell holds a comma separated list like "1:M1,1:M2,..."
#     Replace every bare M<n> token with 5G-M<n>
$reqSheet = $wb.Worksheets.Item("requirements")
$reqUsed = $reqSheet.UsedRange
$reqRowCount = $reqUsed.Rows.Count

for ($r = 2; $r -le $reqRowCount; $r++) {
    $cell = $reqSheet.Cells.Item($r, 6)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $newVal = $val -replace 'M(\d+)', '5G-M$1'
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
